$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column headers "I0" (I1) and "IF" (J1), copying the header style
# (bold, centered, bordered) from the existing H1 header cell.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data for the new I and J columns, row by row (sheet rows 2-69).
$data = @(
    @{Row=2; I=7; J=7},
    @{Row=3; I=7; J=8},
    @{Row=4; I=4; J=6},
    @{Row=5; I=8; J=8},
    @{Row=6; I=5; J=6},
    @{Row=7; I=6; J=6},
    @{Row=8; I=7; J=8},
    @{Row=9; I=7; J=7},
    @{Row=10; I=9; J=9},
    @{Row=11; I=7; J=7},
    @{Row=12; I=7; J=7},
    @{Row=13; I=7; J=7},
    @{Row=14; I=10; J=10},
    @{Row=15; I=8; J=8},
    @{Row=16; I=6; J=7},
    @{Row=17; I=5; J=5},
    @{Row=18; I=6; J=7},
    @{Row=19; I=6; J=7},
    @{Row=20; I=7; J=7},
    @{Row=21; I=5; J=6},
    @{Row=22; I=9; J=9},
    @{Row=23; I=9; J=9},
    @{Row=24; I=7; J=8},
    @{Row=25; I=7; J=7},
    @{Row=26; I=9; J=9},
    @{Row=27; I=7; J=7},
    @{Row=28; I=7; J=8},
    @{Row=29; I=8; J=8},
    @{Row=30; I=7; J=7},
    @{Row=31; I=8; J=9},
    @{Row=32; I=9; J=9},
    @{Row=33; I=8; J=8},
    @{Row=34; I=8; J=8},
    @{Row=35; I=7; J=7},
    @{Row=36; I=9; J=9},
    @{Row=37; I=7; J=8},
    @{Row=38; I=9; J=9},
    @{Row=39; I=10; J=10},
    @{Row=40; I=9; J=9},
    @{Row=41; I=8; J=8},
    @{Row=42; I=9; J=10},
    @{Row=43; I=7; J=7},
    @{Row=44; I=8; J=8},
    @{Row=45; I=9; J=9},
    @{Row=46; I=8; J=8},
    @{Row=47; I=8; J=8},
    @{Row=48; I=8; J=8},
    @{Row=49; I=9; J=9},
    @{Row=50; I=8; J=8},
    @{Row=51; I=7; J=7},
    @{Row=52; I=7; J=7},
    @{Row=53; I=7; J=8},
    @{Row=54; I=8; J=8},
    @{Row=55; I=8; J=8},
    @{Row=56; I=7; J=8},
    @{Row=57; I=9; J=9},
    @{Row=58; I=8; J=9},
    @{Row=59; I=8; J=8},
    @{Row=60; I=9; J=9},
    @{Row=61; I=8; J=8},
    @{Row=62; I=7; J=7},
    @{Row=63; I=8; J=9},
    @{Row=64; I=5; J=5},
    @{Row=65; I=8; J=8},
    @{Row=66; I=5; J=5},
    @{Row=67; I=5; J=5},
    @{Row=68; I=3; J=3},
    @{Row=69; I=4; J=4}
)

foreach ($row in $data) {
    $ws.Cells.Item($row.Row, 9).Value = $row.I
    $ws.Cells.Item($row.Row, 10).Value = $row.J
}
